$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("睡眠日记 Sleep Diary")

# Fill in Day 4 (column E) data for rows 122-135, mirroring the other day columns
$ws.Range("E122").Value = "7：36"
$ws.Range("E123").Value = "8：00"
$ws.Range("E124").Value = "23：10"
$ws.Range("E125").Value = "23：10"
$ws.Range("E126").Value = 5
$ws.Range("E127").Value = 2
$ws.Range("E128").Value = 30
$ws.Range("E129").Value = 490
$ws.Range("E130").Value = "无"
$ws.Range("E131").Value = "无"
$ws.Range("E132").Value = 4
$ws.Range("E133").Value = 3
$ws.Range("E134").Value = 4
$ws.Range("E135").Value = "无"

# Update the view state to match where the user ended up working
$ws.Range("F135").Select()
$excel.ActiveWindow.ScrollRow = 114
